$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to be treated as text so that numeric-looking
# strings (e.g. "1.000", "310.36") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.898.83"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "1.808.96"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "310.36"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").Value = "0.4622"
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("D8").Value = "0.3711"
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "0.07386"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "0.8749"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "20.48"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").Value = "1.845.17"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "5.360"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "92.30"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "6.521"
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").Value = "0.07039"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "0.000008697"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "14.76"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "26.896.35"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "5.323"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("D24").Value = "2.011.47"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("E25").Value = "  -3.46%  "
$ws.Range("D26").Value = "151.35"
$ws.Range("D27").Value = "18.37"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").Value = "2.150"
$ws.Range("E28").Value = "  -6.22%  "
$ws.Range("D29").Value = "5.319"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "115.87"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "0.08916"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "0.7541"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").Value = "1.160"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").Value = "4.448"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D39").Value = "0.05249"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "2.415"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").Value = "2.927"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.5322"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "7.218"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "0.1663"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("D45").Value = "8.523"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("D46").Value = "0.4985"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "10.39"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("D48").Value = "1.675"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").Value = "0.06294"
$ws.Range("E51").Value = "  -1.48%  "
